$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q0)
$ws.Cells.Item(2, 2).Value = -0.0227839658534362
$ws.Cells.Item(2, 3).Value = 0.6421757934121076
$ws.Cells.Item(2, 4).Value = 0.7299546130156783
$ws.Cells.Item(2, 5).Value = 0.8543738133953301
$ws.Cells.Item(2, 6).Value = 0.8732649444381481
$ws.Cells.Item(2, 7).Value = 23

# Row 3 (Q1)
$ws.Cells.Item(3, 2).Value = 0.7639704308143828
$ws.Cells.Item(3, 3).Value = 1.10269471751686
$ws.Cells.Item(3, 4).Value = 2.612846386011487
$ws.Cells.Item(3, 5).Value = 1.616430136446202
$ws.Cells.Item(3, 6).Value = 1.458020551013975
$ws.Cells.Item(3, 7).Value = 22

# Row 4 (Q2)
$ws.Cells.Item(4, 2).Value = 0.6235520873978828
$ws.Cells.Item(4, 3).Value = 1.3761391622945
$ws.Cells.Item(4, 4).Value = 4.053773020290964
$ws.Cells.Item(4, 5).Value = 2.013398375953195
$ws.Cells.Item(4, 6).Value = 1.96168386987362
$ws.Cells.Item(4, 7).Value = 21

# Row 5 (Q3)
$ws.Cells.Item(5, 2).Value = 0.7192758613889139
$ws.Cells.Item(5, 3).Value = 0.8905910875909095
$ws.Cells.Item(5, 4).Value = 1.471357591919539
$ws.Cells.Item(5, 5).Value = 1.212995297566952
$ws.Cells.Item(5, 6).Value = 1.002102960958012
$ws.Cells.Item(5, 7).Value = 20

# Row 6 (Q4)
$ws.Cells.Item(6, 2).Value = 0.541769264226673
$ws.Cells.Item(6, 3).Value = 0.8089314641762144
$ws.Cells.Item(6, 4).Value = 1.078722909040149
$ws.Cells.Item(6, 5).Value = 1.038615862116572
$ws.Cells.Item(6, 6).Value = 0.9104019409704378
$ws.Cells.Item(6, 7).Value = 19

# Row 7 (Q5)
$ws.Cells.Item(7, 2).Value = 0.3378155251717742
$ws.Cells.Item(7, 3).Value = 0.6939632015921933
$ws.Cells.Item(7, 4).Value = 0.7224693225991232
$ws.Cells.Item(7, 5).Value = 0.8499819542785148
$ws.Cells.Item(7, 6).Value = 0.8025803930388514
$ws.Cells.Item(7, 7).Value = 18

# Row 8 (Q6)
$ws.Cells.Item(8, 2).Value = 0.276410312919357
$ws.Cells.Item(8, 3).Value = 0.7156251782359094
$ws.Cells.Item(8, 4).Value = 0.7145716310038928
$ws.Cells.Item(8, 5).Value = 0.8453233884164644
$ws.Cells.Item(8, 6).Value = 0.8234406660685688
$ws.Cells.Item(8, 7).Value = 17

# Row 9 (Q7)
$ws.Cells.Item(9, 2).Value = 0.553324246477921
$ws.Cells.Item(9, 3).Value = 0.722017931190802
$ws.Cells.Item(9, 4).Value = 0.7182827940224992
$ws.Cells.Item(9, 5).Value = 0.8475156600455824
$ws.Cells.Item(9, 6).Value = 0.6705073294552746
$ws.Cells.Item(9, 7).Value = 12

# Row 10 (Q8)
$ws.Cells.Item(10, 2).Value = 0.4056188661620858
$ws.Cells.Item(10, 3).Value = 0.7064691664491064
$ws.Cells.Item(10, 4).Value = 0.7927601511797749
$ws.Cells.Item(10, 5).Value = 0.8903707942086684
$ws.Cells.Item(10, 6).Value = 0.856118606089144
$ws.Cells.Item(10, 7).Value = 7
